$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-23) holds the "Förändrad" date, stored as serial date 45189 (2023-09-20).
# Update it to 45190 (2023-09-21) for every data row, keeping existing formatting.
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
